$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Milestones" worksheet between "Roadmap" and
#    "Technology".
# ------------------------------------------------------------------
$roadmap = $wb.Worksheets.Item("Roadmap")
$technology = $wb.Worksheets.Item("Technology")

$ms = $wb.Worksheets.Add($technology)
$ms.Name = "Milestones"

# ------------------------------------------------------------------
# 2. Populate the header row and the milestone data.
# ------------------------------------------------------------------
$ms.Range("A1").Value = "ID"
$ms.Range("B1").Value = "Date"
$ms.Range("C1").Value = "Name"

$ids = @("MS01","MS02","MS03","MS04","MS05","MS06","MS07","MS08","MS09","MS10")
$names = @(
    "Project Kickoff",
    "Preliminary Design Review (PDR)",
    "Critical Design Review (CDR)",
    "First Prototype Assembled",
    "First Integrated Demo",
    "Software Feature Freeze",
    "Final Demo / Customer Review",
    "Production Decision Point",
    "Initial Manufacturing Run",
    "Product Launch"
)
$dates = @(45703,45748,45792,45848,45870,45915,45945,45962,46037,46082)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ms.Range("A$row").Value = $ids[$i]
    $ms.Range("B$row").Value = $dates[$i]
    $ms.Range("C$row").Value = $names[$i]
}

# ------------------------------------------------------------------
# 3. Re-use the existing cell styles from the "Roadmap" sheet so no
#    new style entries are created in styles.xml.
#      style 5 -> bold header
#      style 7 -> regular text, vertical-centered + wrapped
#      style 8 -> date, vertical-centered + wrapped
# ------------------------------------------------------------------
$roadmap.Range("A1").Copy() | Out-Null
$ms.Range("A1:C1").PasteSpecial(-4122) | Out-Null

$roadmap.Range("A2").Copy() | Out-Null
$ms.Range("A2:A11").PasteSpecial(-4122) | Out-Null
$ms.Range("C2:C11").PasteSpecial(-4122) | Out-Null

$roadmap.Range("C2").Copy() | Out-Null
$ms.Range("B2:B11").PasteSpecial(-4122) | Out-Null

# Left-over formatted (but empty) cells, mirroring the source sheet.
$roadmap.Range("C2").Copy() | Out-Null
$ms.Range("F6:F15").PasteSpecial(-4122) | Out-Null

$roadmap.Range("A2").Copy() | Out-Null
$ms.Range("G6:G15").PasteSpecial(-4122) | Out-Null
$ms.Range("H6:H15").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4. Column widths / layout.
# ------------------------------------------------------------------
$ms.Columns.Item(2).ColumnWidth = 9.5
$ms.Columns.Item(3).ColumnWidth = 33.666666666666664

# ------------------------------------------------------------------
# 5. Selection / activation state.
# ------------------------------------------------------------------
$ms.Activate() | Out-Null
$ms.Range("E13").Select() | Out-Null
